$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "21.689.80"
Set-TextValue "E2" "  -1.31%  "
Set-TextValue "D3" "1.539.77"
Set-TextValue "E3" "  -0.79%  "
Set-TextValue "D4" "1.001"
Set-TextValue "E4" "  -0.01%  "
Set-TextValue "D5" "1.002"
Set-TextValue "D6" "289.97"
Set-TextValue "E6" "  +1.39%  "
Set-TextValue "D7" "0.3910"
Set-TextValue "E7" "  +3.70%  "
Set-TextValue "D8" "0.3171"
Set-TextValue "E8" "  -1.65%  "
Set-TextValue "D9" "42.68"
Set-TextValue "E9" "  +4.24%  "
Set-TextValue "D10" "0.07192"
Set-TextValue "E10" "  -1.16%  "
Set-TextValue "D11" "1.056"
Set-TextValue "E11" "  -5.56%  "
Set-TextValue "D12" "1.002"
Set-TextValue "E12" "  +0.00%  "
Set-TextValue "D13" "5.625"
Set-TextValue "E13" "  -1.23%  "
Set-TextValue "D14" "18.60"
Set-TextValue "E14" "  -3.46%  "
Set-TextValue "D15" "6.614"
Set-TextValue "E15" "  -2.30%  "
Set-TextValue "D16" "1.541.70"
Set-TextValue "E16" "  -0.65%  "
Set-TextValue "D17" "0.00001100"
Set-TextValue "E17" "  +2.22%  "
Set-TextValue "D18" "0.06586"
Set-TextValue "E18" "  -0.66%  "
Set-TextValue "D19" "83.13"
Set-TextValue "E19" "  -2.03%  "
Set-TextValue "E20" "  +0.11%  "
Set-TextValue "D21" "6.148"
Set-TextValue "E21" "  -4.11%  "
Set-TextValue "D22" "15.39"
Set-TextValue "E22" "  -3.15%  "
Set-TextValue "D23" "10.88"
Set-TextValue "E23" "  -4.80%  "
Set-TextValue "D24" "2.390"
Set-TextValue "E24" "  +5.62%  "
Set-TextValue "D25" "21.697.34"
Set-TextValue "E25" "  -1.32%  "
Set-TextValue "D26" "2.350"
Set-TextValue "E26" "  -6.18%  "
Set-TextValue "D27" "147.07"
Set-TextValue "E27" "  -0.74%  "
Set-TextValue "D28" "18.41"
Set-TextValue "E28" "  -1.60%  "
Set-TextValue "D29" "4.837"
Set-TextValue "E29" "  -0.27%  "
Set-TextValue "D30" "1.715.13"
Set-TextValue "E30" "  -0.65%  "
Set-TextValue "D31" "117.40"
Set-TextValue "E31" "  -2.24%  "
Set-TextValue "D32" "0.9681"
Set-TextValue "E32" "  -12.74%  "
Set-TextValue "D33" "5.891"
Set-TextValue "E33" "  -0.07%  "
Set-TextValue "D34" "0.08172"
Set-TextValue "E34" "  +0.62%  "
Set-TextValue "D35" "8.740"
Set-TextValue "E35" "  -5.43%  "
Set-TextValue "D36" "0.06071"
Set-TextValue "E36" "  -1.36%  "
Set-TextValue "D37" "5.115"
Set-TextValue "E37" "  -2.02%  "
Set-TextValue "D38" "0.02198"
Set-TextValue "E38" "  -3.25%  "
Set-TextValue "D39" "0.2029"
Set-TextValue "E39" "  -3.79%  "
Set-TextValue "D40" "1.176"
Set-TextValue "E40" "  -2.63%  "
Set-TextValue "D41" "1.421"
Set-TextValue "E41" "  -14.60%  "
Set-TextValue "D43" "10.68"
Set-TextValue "E43" "  -1.40%  "
Set-TextValue "D44" "0.5727"
Set-TextValue "E44" "  -3.13%  "
Set-TextValue "D45" "3.745"
Set-TextValue "E45" "  +0.92%  "
Set-TextValue "D46" "13.00"
Set-TextValue "E46" "  -3.48%  "
Set-TextValue "D47" "0.5479"
Set-TextValue "E47" "  -3.95%  "
Set-TextValue "D48" "1.160"
Set-TextValue "E48" "  +0.78%  "
Set-TextValue "D49" "116.19"
Set-TextValue "E49" "  -2.64%  "
Set-TextValue "D50" "1.868"
Set-TextValue "E50" "  -2.55%  "
Set-TextValue "D51" "0.06698"
Set-TextValue "E51" "  -2.62%  "
